# Apply the "ix" / zero-indexed coding additions to the Notation workbook.
#
# The "Code" worksheet (sheet2) gets two new rows appended below the
# existing notation table:
#   Row 8: A8 = "i"  (already an existing term), B8 = description of "ix"
#          (the non-code / 1-based index used outside of coding)
#   Row 9: A9 = "ix", B9 = description of the zero-indexed coding convention
#
# The selection/active cell also moves on to B11 (just below the new table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Code")

# Fill column A first (top-to-bottom), then column B (top-to-bottom) so that
# new shared-string entries are created in the same order Excel produced them
# in: "ix", then the "non code index..." text, then the "Zero-indexed..." text.
$ws.Range("A8").Value = "i"
$ws.Range("A9").Value = "ix"

$ws.Range("B8").Value = "non code index e.g. PCA number - starts at 1. This is the most expected index for use outside of zero indexed coding. Used in script to choose PC index etc"
$ws.Range("B9").Value = "Zero-indexed coding - this will delineate between code index and other types of index. Always convert i index input into ix index at start of any functions using it."

# Match the wrap-text style (s="1") used by the rest of the table, and the
# taller 28.8pt row height used by similarly-wrapped rows.
$ws.Range("A8:B9").WrapText = $true
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(9).RowHeight = 28.8

# Move the active selection to B11, as in the saved workbook.
$ws.Range("B11").Select()
